# The source data set gained one more weekly observation (row for
# 2023-12 serial date 45275), which gets inserted right above the
# existing row 558, pushing every subsequent record down by one row
# (old row 558 -> new row 559, ..., old row 675 -> new row 676).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 558; Excel shifts rows 558:675
# down to 559:676 automatically (mirrors Rows("558:558").Insert()).
$ws.Rows(558).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A558").Value = 3
$ws.Range("B558").Value = 'Femacal de La Calera'
$ws.Range("C558").Value = 'Coquimbo'
$ws.Range("D558").Value = 45275
$ws.Range("E558").Value = 5
$ws.Range("F558").Value = 100114013
$ws.Range("G558").Value = 'Zanahoria'
$ws.Range("H558").Value = 'Sin especificar'
$ws.Range("I558").Value = 'Primera'
$ws.Range("J558").Value = 260
$ws.Range("K558").Value = 6000
$ws.Range("L558").Value = 6000
$ws.Range("M558").Value = 6000
$ws.Range("N558").Value = '$/saco 20 kilos'
$ws.Range("O558").Value = 'Provincia de Quillota'
$ws.Range("P558").Value = 300
$ws.Range("Q558").Value = 20
$ws.Range("R558").Value = 'Hortaliza'
